# Applies the diff:
#  - D2 is converted from a text value "26433" to the true number 26433
#  - Five new rows (3-7) are appended, cloning the A/B/C/E values of row 2
#    ("PRIYA", "", "AGARWAL", "NITIN") with varying phone numbers in D:
#       D3 = 26433            (number)
#       D4 = 1234678998765    (number)
#       D5 = 1234678998765    (number)
#       D6 = 23456789         (number)
#       D7 = "2378"           (kept as text, not auto-converted to a number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix D2: store as a real number instead of a numeric-looking text string ---
$ws.Range("D2").Value = 26433

# Reusable values from row 2 (A/B/C/E are identical across rows 2-7)
$firstName = $ws.Range("A2").Value2
$middle    = $ws.Range("B2").Value2
$lastName  = $ws.Range("C2").Value2
$mdlName   = $ws.Range("E2").Value2

$phoneValues = @(26433, 1234678998765, 1234678998765, 23456789)

for ($i = 0; $i -lt $phoneValues.Length; $i++) {
    $r = 3 + $i
    $ws.Range("A$r").Value = $firstName
    $ws.Range("B$r").Value = $middle
    $ws.Range("C$r").Value = $lastName
    $ws.Range("D$r").Value = $phoneValues[$i]
    $ws.Range("E$r").Value = $mdlName
}

# Row 7: D7 must remain text "2378" (not get coerced into the number 2378)
$ws.Range("A7").Value = $firstName
$ws.Range("B7").Value = $middle
$ws.Range("C7").Value = $lastName
$ws.Range("D7").Value = "'2378"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = $mdlName
